$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value (prodNumber changed from ZU706A to 726722-B21)
$ws.Range("B2").Value = "726722-B21"

# Add new columns H and I with headers "user"/"pass" and values "test"/"apassword"
$ws.Range("H1").Value = "user"
$ws.Range("I1").Value = "pass"
$ws.Range("H2").Value = "test"
$ws.Range("I2").Value = "apassword"

# Match the style used by the other data cells in row 2 (copy formats from C2,
# which already carries the shared "data row" style, onto the cells that need it)
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to match the diff (activeCell J9)
$ws.Range("J9").Select()
